$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / recalculated mean
$ws.Range("F2").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 0
$ws.Range("F15").Value = -1
$ws.Range("F19").Value = 0
